$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 198.07692
$ws.Range("I9").Value = 59.75
$ws.Range("K9").Value = 59.75
$ws.Range("M9").Value = 109.25

$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()

$ws.Range("H40").Value = 13014.714
$ws.Range("J40").Value = 14133.833
$ws.Range("L40").Value = 14133.833
$ws.Range("N40").Value = -14483.833

$ws.Range("H43").Value = 5822.5
$ws.Range("I43").Value = 6876.778
$ws.Range("J43").Value = 4241.0835
$ws.Range("K43").Value = 6876.778
$ws.Range("L43").Value = 4241.0835
$ws.Range("M43").Value = -6807.778
$ws.Range("N43").Value = -4379.0835

$ws.Range("H64").Value = 4922.4707
$ws.Range("J64").Value = 4323.625
$ws.Range("L64").Value = 4323.625
$ws.Range("N64").Value = -4819.625

$ws.Range("H67").Value = 4922.4707
$ws.Range("J67").Value = 4323.625
$ws.Range("L67").Value = 4323.625
$ws.Range("N67").Value = -6039.625

$ws.Range("H80").Value = 526.9231
$ws.Range("I80").Value = 382.77777
$ws.Range("K80").Value = 1148.33331
$ws.Range("M80").Value = -150.33331

$ws.Range("H83").Value = 526.9231
$ws.Range("I83").Value = 382.77777
$ws.Range("K83").Value = 3444.99993
$ws.Range("M83").Value = 1547.00007

$ws.Range("H129").Value = 4179.5
$ws.Range("J129").Value = 2077
$ws.Range("L129").Value = 6231
$ws.Range("N129").Value = -16231

$ws.Range("H132").Value = 2625.0576
$ws.Range("J132").Value = 2998.25
$ws.Range("L132").Value = 8994.75
$ws.Range("N132").Value = -14054.75

$ws.Range("H134").Value = 60980.77
$ws.Range("J134").Value = 60980.77
$ws.Range("L134").Value = 60980.77
$ws.Range("N134").Value = -71120.76999999999

$ws.Range("H138").Value = 3010.288
$ws.Range("I138").Value = 2349.4375
$ws.Range("K138").Value = 7048.3125
$ws.Range("M138").Value = -1908.3125

$ws.Range("H141").Value = 3587.6
$ws.Range("I141").Value = 1687.6296
$ws.Range("K141").Value = 5062.8888
$ws.Range("M141").Value = 117.1112000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3965.672
$ws.Range("I32").Value = 2944.818
$ws.Range("K32").Value = 2944.818
$ws.Range("M32").Value = -2657.818

$ws.Range("H74").Value = 3590.6667
$ws.Range("I74").Value = 1512.9
$ws.Range("K74").Value = 1512.9
$ws.Range("M74").Value = -638.9000000000001

$ws.Range("H77").Value = 3590.6667
$ws.Range("I77").Value = 1512.9
$ws.Range("K77").Value = 7564.5
$ws.Range("M77").Value = -3196.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2985.2778
$ws.Range("I20").Value = 3011.1155
$ws.Range("J20").Value = 2918.1
$ws.Range("K20").Value = 3011.1155
$ws.Range("L20").Value = 2918.1
$ws.Range("M20").Value = -2764.1155
$ws.Range("N20").Value = -3412.1

$ws.Range("H27").Value = 39995
$ws.Range("J27").Value = 39995
$ws.Range("L27").Value = 39995
$ws.Range("N27").Value = -40379

$ws.Range("H75").Value = 121762.22
$ws.Range("I75").Value = 12070
$ws.Range("K75").Value = 12070
$ws.Range("M75").Value = -11134

$ws.Range("H78").Value = 121762.22
$ws.Range("I78").Value = 12070
$ws.Range("K78").Value = 36210
$ws.Range("M78").Value = -31530

$ws.Range("H97").Value = 94617.664
$ws.Range("I97").Value = 5960.3335
$ws.Range("K97").Value = 5960.3335
$ws.Range("M97").Value = -4969.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = 12692.714
$ws.Range("I39").Value = 8141.6665
$ws.Range("K39").Value = 8141.6665
$ws.Range("M39").Value = -7750.6665

$ws.Range("H49").Value = 12692.714
$ws.Range("I49").Value = 8141.6665
$ws.Range("K49").Value = 8141.6665
$ws.Range("M49").Value = -7959.6665

$ws.Range("H94").Value = 2249.0588
$ws.Range("I94").Value = 1951.3636
$ws.Range("K94").Value = 1951.3636
$ws.Range("M94").Value = -1500.3636

$ws.Range("H132").Value = 4642.706
$ws.Range("I132").Value = 3592.182
$ws.Range("K132").Value = 10776.546
$ws.Range("M132").Value = -8246.545999999998

$ws.Range("H134").Value = 4695.7334
$ws.Range("I134").Value = 3647.2173
$ws.Range("K134").Value = 10941.6519
$ws.Range("M134").Value = -8406.651899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 86329.664
$ws.Range("J37").Value = 86329.664
$ws.Range("L37").Value = 258988.992
$ws.Range("N37").Value = -259212.992

$ws.Range("H86").Value = 9539.166999999999
$ws.Range("I86").Value = 11166.4
$ws.Range("K86").Value = 33499.2
$ws.Range("M86").Value = -32313.2

$ws.Range("H89").Value = 9539.166999999999
$ws.Range("I89").Value = 11166.4
$ws.Range("K89").Value = 100497.6
$ws.Range("M89").Value = -94569.59999999999

$ws.Range("H114").Value = 1142.7142
$ws.Range("I114").Value = 417
$ws.Range("J114").Value = 1868.4286
$ws.Range("K114").Value = 1251
$ws.Range("L114").Value = 5605.2858
$ws.Range("M114").Value = 2003
$ws.Range("N114").Value = -12113.2858

$ws.Range("H129").Value = 11113514
$ws.Range("J129").Value = 16670191
$ws.Range("L129").Value = 50010573
$ws.Range("N129").Value = -50020573

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8465.333000000001
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 5392.5
$ws.Range("I132").Value = 4060.35
$ws.Range("J132").Value = 9833
$ws.Range("K132").Value = 12181.05
$ws.Range("L132").Value = 29499
$ws.Range("M132").Value = -9651.049999999999
$ws.Range("N132").Value = -34559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1533.3334
$ws.Range("I22").Value = 942.8570999999999
$ws.Range("J22").Value = 2050
$ws.Range("K22").Value = 942.8570999999999
$ws.Range("L22").Value = 2050
$ws.Range("M22").Value = -647.8570999999999
$ws.Range("N22").Value = -2640

$ws.Range("H27").Value = 1533.3334
$ws.Range("I27").Value = 942.8570999999999
$ws.Range("J27").Value = 2050
$ws.Range("K27").Value = 942.8570999999999
$ws.Range("L27").Value = 2050
$ws.Range("M27").Value = -835.8570999999999
$ws.Range("N27").Value = -2264

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 63640.25
$ws.Range("I81").Value = 84407.25
$ws.Range("J81").Value = 1339.25
$ws.Range("K81").Value = 168814.5
$ws.Range("L81").Value = 2678.5
$ws.Range("M81").Value = -167753.5
$ws.Range("N81").Value = -4800.5

$ws.Range("H84").Value = 63640.25
$ws.Range("I84").Value = 84407.25
$ws.Range("J84").Value = 1339.25
$ws.Range("K84").Value = 844072.5
$ws.Range("L84").Value = 13392.5
$ws.Range("M84").Value = -838768.5
$ws.Range("N84").Value = -24000.5

$ws.Range("H122").Value = 3750.389
$ws.Range("I122").Value = 3509.7856
$ws.Range("K122").Value = 10529.3568
$ws.Range("M122").Value = -8079.356800000001
